# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Actualizar "Periodo Mora" (columna E, filas 16-22) a orden ascendente ---
#        (antes: 2207,2206,2205,2204,2203,2202,2201 -> ahora: 2201..2207)
$periodos = @("2201", "2202", "2203", "2204", "2205", "2206", "2207")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $fila = 16 + $i
    $ws.Range("E" + $fila).Value = $periodos[$i]
}

# --- 2. Actualizar "Salario Basico" (columna G, filas 16-22): 908526 -> 877803 ---
for ($fila = 16; $fila -le 22; $fila++) {
    $ws.Range("G" + $fila).Value = 877803
}

# --- 3. Mover el logo (imagen) 19pt hacia la izquierda ---
$shp = $ws.Shapes.Item(1)
$shp.Left = $shp.Left - 19

# --- 4. Ajustar anchos de columnas B:J (recalculo de autofit tras la actualizacion) ---
$ws.Columns.Item(2).ColumnWidth = 16.07
$ws.Columns.Item(3).ColumnWidth = 7.98
$ws.Columns.Item(4).ColumnWidth = 26.53
$ws.Columns.Item(5).ColumnWidth = 11.89
$ws.Columns.Item(6).ColumnWidth = 8.62
$ws.Columns.Item(7).ColumnWidth = 12.62
$ws.Columns.Item(8).ColumnWidth = 17.07
$ws.Columns.Item(9).ColumnWidth = 15.98
$ws.Columns.Item(10).ColumnWidth = 13.35

Write-Host "Cambios aplicados correctamente"
